# Apply "updated example data per email from bjorn and tests passing":
#  - add a new "Renames" worksheet (old name / new name lookup table) as the
#    last sheet in the workbook, and make it the active sheet/tab
#  - clear the stale selection that was left on "Composition IDs"

$wb = $excel.ActiveWorkbook

$compositionIds = $wb.Worksheets.Item("Composition IDs")

# Insert the new sheet after the last existing sheet so it lands at the end
# of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$renames = $wb.Worksheets.Add($null, $lastSheet)
$renames.Name = "Renames"

# Populate header row first (B1 then A1) and the two data rows so that the
# strings get appended to the shared-string table in "new name" / "old name"
# order, matching how the workbook was authored.
$renames.Range("B1").Value = "new name"
$renames.Range("A1").Value = "old name"

$renames.Range("A2").Value = "Cylinder"
$renames.Range("B2").Value = "Big Cylinder"

$renames.Range("A3").Value = "Lug Nut"
$renames.Range("B3").Value = "Locking Nut"

# Size the columns to fit their contents like the source workbook does.
$renames.Columns("A:A").AutoFit() | Out-Null
$renames.Columns("B:B").AutoFit() | Out-Null

# Leave the cursor below the table and make this the selected/active tab.
$renames.Range("A4").Select()

# The "Composition IDs" sheet no longer keeps the tab-selected/old selection
# state now that "Renames" is the active tab; reset its lingering selection
# back to the top-left cell.
$compositionIds.Range("A1").Select()
$renames.Activate()
